$d = $word.ActiveDocument

# The paragraph currently stores the tag "<id>p123r_1</id>" split across three
# runs: "<id>" (Courier New / color 7f6000), "p123r_1" (plain black run) and
# "</id>" (Courier New / color 7f6000). Collapse them into a single run that
# carries the whole literal string, keeping the formatting of the run that
# Find/Replace seeds from (the first, "<id>", run) for the merged text.
$rng = $d.Content
$rng.Find.Execute("<id>p123r_1</id>", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "<id>p123r_1</id>", 2)
